# Financial Statement: add "wealth class" rows (MMN, LAE, PPAL, Fixed Asset)
# to both the Wealth Allocation block (rows 12-15) and the Cash Flow block
# (rows 18-21), inserting two fresh rows for the first block and one fresh
# row at the end for the second block.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Financial Statement")

# --- Insert two new rows right after row 12 (old rows 13.. shift down by 2) ---
$ws.Range("A13:A14").EntireRow.Insert()

# --- Apply a full thin box border ("boxed" look) to the two new rows ---
$newRows = $ws.Range("A13:N14")
$newRows.Borders.LineStyle = 1
$newRows.Borders.Weight = 2

# --- Fill in the wealth-class labels for the first (Wealth Allocation) block ---
$ws.Range("A12").Value = "MMN"
$ws.Range("A13").Value = "LAE"
$ws.Range("A14").Value = "PPAL"
$ws.Range("A15").Value = "Fixed Asset"

# --- Fill in the wealth-class labels for the second (Cash Flow) block ---
$ws.Range("A18").Value = "MMN"
$ws.Range("A19").Value = "LAE"
$ws.Range("A20").Value = "PPAL"

# --- Append a brand-new row 21 mirroring the new-row formatting/content ---
$lastRow = $ws.Range("A21:N21")
$lastRow.Borders.LineStyle = 1
$lastRow.Borders.Weight = 2
$ws.Range("A21").Value = "Fixed Asset"
Write-Host "Edit complete"
